$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "AddCustomerTest"

# Columns A-C first (firstname / lastname / postcode data)
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

$ws.Range("A2").Value = "sandip"
$ws.Range("B2").Value = "thopate"
$ws.Range("C2").Value = "hsf894r"

$ws.Range("A3").Value = "lasdjf"
$ws.Range("B3").Value = "lslsjf"
$ws.Range("C3").Value = 567898

# Column D last (alerttext data)
$ws.Range("D1").Value = "alerttext"
$ws.Range("D2").Value = "Customer added successfully"
$ws.Range("D3").Value = "Customer added successfully"

# Left alignment for C2:C3
$ws.Range("C2:C3").HorizontalAlignment = -4131

# Column D width (target stored width is 33.109375; the host quantizes
# ColumnWidth to 1/6-character increments, so 32.3 lands on the nearest
# representable stored width of 33.1666...)
$ws.Columns.Item(4).ColumnWidth = 32.3

# Selection
$ws.Range("D3").Select()
